$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (single dot) need an
# explicit text NumberFormat first, otherwise Excel auto-converts the
# assigned string into a numeric value (losing the trailing-zero / exact
# text representation used by this price-ticker sheet).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

$ws.Range("D2").Value = "62.093.91"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").Value = "3.421.61"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("D5").Value = "578.42"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("E6").Value = "  +3.32%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").Value = "8.05"
$ws.Range("E9").Value = "  +4.82%  "
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("E11").Value = "  +3.06%  "
$ws.Range("D12").Value = "4.008.01"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "28.72"
$ws.Range("E14").Value = "  -3.35%  "
$ws.Range("D15").Value = "3.417.81"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "62.150.57"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "6.50"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").Value = "14.56"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").Value = "8.95"
$ws.Range("E20").Value = "  -4.29%  "
$ws.Range("D21").Value = "382.78"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").Value = "0.570"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "75.19"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "3.562.60"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("E26").Value = "  -3.58%  "
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").Value = "7.68"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.12"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "7.91"
$ws.Range("E31").Value = "  -4.25%  "
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").Value = "23.24"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").Value = "5.47"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "168.63"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "31.11"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("D40").Value = "3.455.87"
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("D42").Value = "42.76"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").Value = "4.40"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("D46").Value = "1.17"
$ws.Range("E46").Value = "  -3.04%  "
$ws.Range("D47").Value = "2.542.03"
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E49").Value = "  -4.96%  "
$ws.Range("E51").Value = "  +0.03%  "
